# Updated cryptos list with latest price and volume(1h) data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text like "76.503.76" or "1.00" that must remain text,
# not be auto-converted to numbers, so force the whole Price column to Text format first.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "76.503.76"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.053.32"
$ws.Range("E3").Value = "  +4.60%  "
$ws.Range("D5").Value = "202.17"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "625.26"
$ws.Range("E6").Value = "  +4.75%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +5.13%  "
$ws.Range("D10").Value = "3.051.72"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("D11").Value = "0.441"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +5.64%  "
$ws.Range("D14").Value = "3.616.95"
$ws.Range("E14").Value = "  +4.64%  "
$ws.Range("D15").Value = "29.48"
$ws.Range("D16").Value = "76.386.70"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "3.058.68"
$ws.Range("E18").Value = "  +5.19%  "
$ws.Range("D19").Value = "13.60"
$ws.Range("E19").Value = "  +4.82%  "
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("D21").Value = "375.67"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "4.37"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "73.64"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").Value = "3.209.76"
$ws.Range("E25").Value = "  +4.60%  "
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "9.91"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("E29").Value = "  +3.38%  "
$ws.Range("D30").Value = "0.995"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("E31").Value = "  +7.24%  "
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").Value = "506.89"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +6.80%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "20.89"
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("D37").Value = "162.70"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").Value = "0.389"
$ws.Range("E38").Value = "  +5.55%  "
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "192.92"
$ws.Range("E40").Value = "  +5.29%  "
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +3.87%  "
$ws.Range("D45").Value = "0.807"
$ws.Range("E45").Value = "  +22.76%  "
$ws.Range("D46").Value = "1.28"
$ws.Range("E46").Value = "  +7.54%  "
$ws.Range("D47").Value = "42.18"
$ws.Range("E47").Value = "  +5.09%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "2.47"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("E50").Value = "  +7.16%  "
$ws.Range("E51").Value = "  +6.68%  "
